$d = $word.ActiveDocument

# Word table-cell ranges end with a cell-mark (CR + BEL, chars 13/7); strip
# those before comparing cell text against plain strings.
function CellText($cell) {
  return $cell.Range.Text.TrimEnd([char]13, [char]7)
}

# Locate the "Transformed Data Set" fairness-metric table: a 5-column table
# whose header row reads Mental Health (Race) / Mental Health (Gender) /
# Physical Health (Race) / Physical Health (Gender), and whose
# "Mental Health" cells for the Demographic Parity / Equal accuracy rows
# are still blank (the companion "Original Data Set" table above it
# already has its values filled in).
$targetTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
  $t = $d.Tables.Item($i)
  if ($t.Columns.Count -ne 5 -or $t.Rows.Count -lt 3) { continue }
  $header2 = CellText $t.Cell(1,2)
  $header3 = CellText $t.Cell(1,3)
  if ($header2 -eq "Mental Health (Race)" -and $header3 -eq "Mental Health (Gender)") {
    $c22 = CellText $t.Cell(2,2)
    $c23 = CellText $t.Cell(2,3)
    if ($c22 -eq "" -and $c23 -eq "") {
      $targetTable = $t
      break
    }
  }
}

if ($targetTable -ne $null) {
  $targetTable.Cell(2,2).Range.Text = "1"
  $targetTable.Cell(2,3).Range.Text = "1"
  $targetTable.Cell(3,2).Range.Text = "1"
  $targetTable.Cell(3,3).Range.Text = "1"
  Write-Output "Updated fairness-metric table cells."
} else {
  Write-Output "Target table not found."
}
